## Applies the January-26 simulator-run edits:
##  - Employee ID regenerated
##  - Client names replaced for the Jan26-Jan30 rows
##  - Hours normalized to a flat 8/day (was 9/9/9/9/4+5 split reg/OT)
##  - Rate now 100, Total now 800 per day (was 0/0)
##  - The separate OT line for 2026-01-30 is removed (full-month/day coverage
##    is now captured by the single Regular row, so OT drops to 0)
##  - Subtotal hours 40 (down from 45), OT text updated, HOURLY/GRAND totals recomputed

$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: "Weekly Timesheet" ----------
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# New client names for each work day
$ws1.Range("B2").Value = "Markfield"
$ws1.Range("B3").Value = "Leixner/Smith"
$ws1.Range("B4").Value = "Hunter"
$ws1.Range("B5").Value = "O'Connor"
$ws1.Range("B6").Value = "Varricchio"

# Flat 8 hours/day, $100 rate, $800 total for every day
$ws1.Range("C2").Value = 8
$ws1.Range("C3").Value = 8
$ws1.Range("C4").Value = 8
$ws1.Range("C5").Value = 8
$ws1.Range("C6").Value = 8

$ws1.Range("E2:E6").Value = 100
$ws1.Range("F2:F6").Value = 800

# Remove the separate OT row for 2026-01-30 (row 7) - shifts rows 9-14 up to 8-13
$ws1.Rows.Item(7).Delete()

# Subtotal row (now row 8 after the delete)
$ws1.Range("C8").Value = 40
$ws1.Range("D8").Value = "Reg: 40 / OT: 0"
$ws1.Range("F8").Value = 4000

# HOURLY SUBTOTAL (now row 11) and GRAND TOTAL (now row 13)
$ws1.Range("F11").Value = 4000
$ws1.Range("F13").Value = 4000

# ---------- Sheet 2: "Jason Schema" ----------
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Employee ID regenerated for every row
$ws2.Range("B2").Value = "emp_qhpjptqm"
$ws2.Range("B3").Value = "emp_qhpjptqm"
$ws2.Range("B4").Value = "emp_qhpjptqm"
$ws2.Range("B5").Value = "emp_qhpjptqm"
$ws2.Range("B6").Value = "emp_qhpjptqm"

# Mirror the client-name changes
$ws2.Range("D2").Value = "Markfield"
$ws2.Range("D3").Value = "Leixner/Smith"
$ws2.Range("D4").Value = "Hunter"
$ws2.Range("D5").Value = "O'Connor"
$ws2.Range("D6").Value = "Varricchio"

# Mirror the flat 8 hours/day, $100 rate, $800 total
$ws2.Range("E2").Value = 8
$ws2.Range("E3").Value = 8
$ws2.Range("E4").Value = 8
$ws2.Range("E5").Value = 8
$ws2.Range("E6").Value = 8

$ws2.Range("F2:F6").Value = 100
$ws2.Range("G2:G6").Value = 800

# Remove the separate OT row for 2026-01-30 (row 7)
$ws2.Rows.Item(7).Delete()
